$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.226.51'
$ws.Range("E2").Value = '  +3.30%  '
$ws.Range("D3").Value = '1.899.47'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").Value = "'325.94"
$ws.Range("E5").Value = '  +3.47%  '
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("D7").Value = "'0.5183"
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("D8").Value = "'0.4019"
$ws.Range("E8").Value = '  +2.41%  '
$ws.Range("D9").Value = "'0.08447"
$ws.Range("E9").Value = '  +0.34%  '
$ws.Range("D10").Value = "'42.77"
$ws.Range("E10").Value = '  +0.71%  '
$ws.Range("E11").Value = '  +0.38%  '
$ws.Range("D12").Value = "'23.34"
$ws.Range("E12").Value = '  +12.98%  '
$ws.Range("D13").Value = "'6.440"
$ws.Range("E13").Value = '  +3.20%  '
$ws.Range("D14").Value = '1.898.85'
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("D15").Value = "'7.335"
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("D17").Value = "'94.80"
$ws.Range("E17").Value = '  +1.82%  '
$ws.Range("D19").Value = "'0.06668"
$ws.Range("E19").Value = '  -1.21%  '
$ws.Range("D20").Value = "'18.29"
$ws.Range("E20").Value = '  +2.55%  '
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").Value = "'5.958"
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("D23").Value = '30.232.90'
$ws.Range("E23").Value = '  +3.23%  '
$ws.Range("D24").Value = "'11.28"
$ws.Range("E24").Value = '  +1.24%  '
$ws.Range("D25").Value = "'2.217"
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '2.121.21'
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").Value = "'21.77"
$ws.Range("D28").Value = "'161.35"
$ws.Range("E28").Value = '  +1.38%  '
$ws.Range("D29").Value = "'2.395"
$ws.Range("E29").Value = '  -1.58%  '
$ws.Range("D30").Value = "'129.74"
$ws.Range("E30").Value = '  +1.35%  '
$ws.Range("D31").Value = "'1.092"
$ws.Range("E31").Value = '  +3.21%  '
$ws.Range("D32").Value = "'0.1059"
$ws.Range("E32").Value = '  +1.08%  '
$ws.Range("D33").Value = "'6.026"
$ws.Range("E33").Value = '  -1.63%  '
$ws.Range("D34").Value = "'3.708"
$ws.Range("E34").Value = '  +1.70%  '
$ws.Range("D35").Value = "'0.02493"
$ws.Range("E35").Value = '  +0.66%  '
$ws.Range("D36").Value = "'0.06571"
$ws.Range("E36").Value = '  +0.49%  '
$ws.Range("D37").Value = "'0.2210"
$ws.Range("E37").Value = '  +0.75%  '
$ws.Range("D38").Value = "'5.244"
$ws.Range("E38").Value = '  +2.31%  '
$ws.Range("E39").Value = '  -0.67%  '
$ws.Range("D40").Value = "'11.84"
$ws.Range("E40").Value = '  +5.14%  '
$ws.Range("D41").Value = "'8.784"
$ws.Range("E41").Value = '  -2.87%  '
$ws.Range("D42").Value = "'0.6507"
$ws.Range("D43").Value = "'1.236"
$ws.Range("E43").Value = '  +0.27%  '
$ws.Range("D44").Value = "'0.6106"
$ws.Range("E44").Value = '  +0.77%  '
$ws.Range("D45").Value = "'13.19"
$ws.Range("E45").Value = '  -0.19%  '
$ws.Range("D46").Value = "'3.712"
$ws.Range("D47").Value = "'2.057"
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("D49").Value = "'124.80"
$ws.Range("E49").Value = '  +1.19%  '
$ws.Range("D50").Value = "'1.165"
$ws.Range("E50").Value = '  -0.90%  '
$ws.Range("D51").Value = "'79.14"
$ws.Range("E51").Value = '  +1.93%  '
